$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058280065418936
$ws.Range("D2").Value = 1.064895205961293
$ws.Range("E2").Value = 1.0715832656905
$ws.Range("F2").Value = 1.078129124377686
$ws.Range("I2").Value = 1.051649635205988
$ws.Range("J2").Value = 1.063272107644488
$ws.Range("K2").Value = 1.06760990966836
$ws.Range("L2").Value = 1.074280052769532
$ws.Range("M2").Value = 1.080808608342185
$ws.Range("N2").Value = 1.064782075996358
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059221215800529
$ws.Range("D3").Value = 1.065651726856674
$ws.Range("E3").Value = 1.072475833959087
$ws.Range("F3").Value = 1.079019839379666
$ws.Range("I3").Value = 1.051906294985497
$ws.Range("J3").Value = 1.063866100710597
$ws.Range("K3").Value = 1.068181871349442
$ws.Range("L3").Value = 1.074989014390848
$ws.Range("M3").Value = 1.081516968823583
$ws.Range("N3").Value = 1.065376912600753
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059830865535453
$ws.Range("D4").Value = 1.066141817454895
$ws.Range("E4").Value = 1.073054350186506
$ws.Range("F4").Value = 1.079597122554772
$ws.Range("I4").Value = 1.052071553652995
$ws.Range("J4").Value = 1.064250455460134
$ws.Range("K4").Value = 1.06855187261569
$ws.Range("L4").Value = 1.075448083754015
$ws.Range("M4").Value = 1.081975614049876
$ws.Range("N4").Value = 1.065761813178121
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.060087319666385
$ws.Range("D5").Value = 1.066347986744831
$ws.Range("E5").Value = 1.073297787858388
$ws.Range("F5").Value = 1.07984003370735
$ws.Range("I5").Value = 1.052140832098165
$ws.Range("J5").Value = 1.064412037470876
$ws.Range("K5").Value = 1.068707396993616
$ws.Range("L5").Value = 1.075641152749419
$ws.Range("M5").Value = 1.082168496288153
$ws.Range("N5").Value = 1.065923624653852
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060130388636979
$ws.Range("D6").Value = 1.066382611360722
$ws.Range("E6").Value = 1.073338675524924
$ws.Range("F6").Value = 1.079880832490332
$ws.Range("I6").Value = 1.052152452724993
$ws.Range("J6").Value = 1.064439167727298
$ws.Range("K6").Value = 1.06873350877187
$ws.Range("L6").Value = 1.075673574316591
$ws.Range("M6").Value = 1.082200886001515
$ws.Range("N6").Value = 1.065950793438349
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059834291671217
$ws.Range("D7").Value = 1.066144571769658
$ws.Range("E7").Value = 1.073057602112942
$ws.Range("F7").Value = 1.079600367477878
$ws.Range("I7").Value = 1.052072480127056
$ws.Range("J7").Value = 1.064252614531083
$ws.Range("K7").Value = 1.068553950836359
$ws.Range("L7").Value = 1.07545066325343
$ws.Range("M7").Value = 1.081978191087077
$ws.Range("N7").Value = 1.065763975315198
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058597993684032
$ws.Range("D8").Value = 1.065150756613654
$ws.Range("E8").Value = 1.071884712604269
$ws.Range("F8").Value = 1.078429951961195
$ws.Range("I8").Value = 1.051736543450046
$ws.Range("J8").Value = 1.063472849412765
$ws.Range("K8").Value = 1.067803226034447
$ws.Range("L8").Value = 1.074519581856131
$ws.Range("M8").Value = 1.081047941530954
$ws.Range("N8").Value = 1.064983102840972
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.056424606821709
$ws.Range("D9").Value = 1.063403966509844
$ws.Range("E9").Value = 1.069825389451866
$ws.Range("F9").Value = 1.076374730560487
$ws.Range("I9").Value = 1.051138347406468
$ws.Range("J9").Value = 1.062098863443554
$ws.Range("K9").Value = 1.06647966619317
$ws.Range("L9").Value = 1.072881433930723
$ws.Range("M9").Value = 1.079410991591309
$ws.Range("N9").Value = 1.063607165654087
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.05497920233662
$ws.Range("D10").Value = 1.062242513195333
$ws.Range("E10").Value = 1.06845761329232
$ws.Range("F10").Value = 1.075009517832128
$ws.Range("I10").Value = 1.050735396275777
$ws.Range("J10").Value = 1.061182977574536
$ws.Range("K10").Value = 1.065596898071863
$ws.Range("L10").Value = 1.071791116644473
$ws.Range("M10").Value = 1.078321294872038
$ws.Range("N10").Value = 1.062689979122081
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.054354176117638
$ws.Range("D11").Value = 1.061740338688256
$ws.Range("E11").Value = 1.067866581361724
$ws.Range("F11").Value = 1.074419555796281
$ws.Range("I11").Value = 1.050559936438386
$ws.Range("J11").Value = 1.060786427218124
$ws.Range("K11").Value = 1.065214570055811
$ws.Range("L11").Value = 1.071319435617077
$ws.Range("M11").Value = 1.077849840740717
$ws.Range("N11").Value = 1.062292865618678
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054122141442156
$ws.Range("D12").Value = 1.061553921751551
$ws.Range("E12").Value = 1.067647231239413
$ws.Range("F12").Value = 1.074200597096852
$ws.Range("I12").Value = 1.050494616302898
$ws.Range("J12").Value = 1.060639136973045
$ws.Range("K12").Value = 1.065072544955579
$ws.Range("L12").Value = 1.071144298680205
$ws.Range("M12").Value = 1.077674781835775
$ws.Range("N12").Value = 1.062145366204557
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054171907861482
$ws.Range("D13").Value = 1.061593903678283
$ws.Range("E13").Value = 1.067694274168482
$ws.Range("F13").Value = 1.07424755633453
$ws.Range("I13").Value = 1.050508634315332
$ws.Range("J13").Value = 1.060670730932291
$ws.Range("K13").Value = 1.065103010322117
$ws.Range("L13").Value = 1.071181863122731
$ws.Range("M13").Value = 1.077712329822339
$ws.Range("N13").Value = 1.062177005030848
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.054334993443236
$ws.Range("D14").Value = 1.061724927091619
$ws.Range("E14").Value = 1.067848446020134
$ws.Range("F14").Value = 1.074401452931666
$ws.Range("I14").Value = 1.050554540041519
$ws.Range("J14").Value = 1.060774252027591
$ws.Range("K14").Value = 1.065202830443924
$ws.Range("L14").Value = 1.071304957381429
$ws.Range("M14").Value = 1.077835369082457
$ws.Range("N14").Value = 1.062280673137978
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.054435492834668
$ws.Range("D15").Value = 1.061805669963682
$ws.Range("E15").Value = 1.067943461006428
$ws.Range("F15").Value = 1.074496297532843
$ws.Range("I15").Value = 1.050582804674449
$ws.Range("J15").Value = 1.060838035646319
$ws.Range("K15").Value = 1.065264331431121
$ws.Range("L15").Value = 1.071380808653048
$ws.Range("M15").Value = 1.077911185640833
$ws.Range("N15").Value = 1.06234454733676
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.055020701138871
$ws.Range("D16").Value = 1.06227585663315
$ws.Range("E16").Value = 1.068496864057861
$ws.Range("F16").Value = 1.075048696751012
$ws.Range("I16").Value = 1.050747020384499
$ws.Range("J16").Value = 1.06120929612855
$ws.Range("K16").Value = 1.065622270256338
$ws.Range("L16").Value = 1.071822429802929
$ws.Range("M16").Value = 1.078352592094144
$ws.Range("N16").Value = 1.06271633505146
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.055388013490382
$ws.Range("D17").Value = 1.062570992133915
$ws.Range("E17").Value = 1.068844327999502
$ws.Range("F17").Value = 1.075395520309539
$ws.Range("I17").Value = 1.050849766783187
$ws.Range("J17").Value = 1.061442187946347
$ws.Range("K17").Value = 1.065846774225096
$ws.Range("L17").Value = 1.072099564047165
$ws.Range("M17").Value = 1.078629580494697
$ws.Range("N17").Value = 1.062949557602351
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.055602341873504
$ws.Range("D18").Value = 1.062743211213117
$ws.Range("E18").Value = 1.069047115849912
$ws.Range("F18").Value = 1.075597930764801
$ws.Range("I18").Value = 1.050909602487147
$ws.Range("J18").Value = 1.061578033038311
$ws.Range("K18").Value = 1.065977715451762
$ws.Range("L18").Value = 1.072261253509855
$ws.Range("M18").Value = 1.078791180852516
$ws.Range("N18").Value = 1.063085595609927
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.055675436068256
$ws.Range("D19").Value = 1.062801945556799
$ws.Range("E19").Value = 1.069116281236606
$ws.Range("F19").Value = 1.075666966811806
$ws.Range("I19").Value = 1.050929988863569
$ws.Range("J19").Value = 1.06162435322044
$ws.Range("K19").Value = 1.066022361611896
$ws.Range("L19").Value = 1.072316392480823
$ws.Range("M19").Value = 1.078846288755274
$ws.Range("N19").Value = 1.063131981572028
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05534859593326
$ws.Range("D20").Value = 1.062539319477858
$ws.Range("E20").Value = 1.068807036189323
$ws.Range("F20").Value = 1.075358297606472
$ws.Range("I20").Value = 1.050838752840588
$ws.Range("J20").Value = 1.061417200516068
$ws.Range("K20").Value = 1.065822687907892
$ws.Range("L20").Value = 1.072069825838348
$ws.Range("M20").Value = 1.078599858350815
$ws.Range("N20").Value = 1.062924534687055
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.054286965264342
$ws.Range("D21").Value = 1.061686340829778
$ws.Range("E21").Value = 1.067803041120726
$ws.Range("F21").Value = 1.074356129252719
$ws.Range("I21").Value = 1.050541025989741
$ws.Range("J21").Value = 1.060743767465123
$ws.Range("K21").Value = 1.065173436213228
$ws.Range("L21").Value = 1.07126870734722
$ws.Range("M21").Value = 1.077799135416139
$ws.Range("N21").Value = 1.062250145283935
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.053620216452081
$ws.Range("D22").Value = 1.061150694020858
$ws.Range("E22").Value = 1.067172863304253
$ws.Range("F22").Value = 1.073727065363978
$ws.Range("I22").Value = 1.050352985766181
$ws.Range("J22").Value = 1.060320389857331
$ws.Range("K22").Value = 1.064765160196656
$ws.Range("L22").Value = 1.070765397142188
$ws.Range("M22").Value = 1.077296037784854
$ws.Range("N22").Value = 1.061826166431377
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.053973602026852
$ws.Range("D23").Value = 1.06143458797088
$ws.Range("E23").Value = 1.067506830197847
$ws.Range("F23").Value = 1.074060445010109
$ws.Range("I23").Value = 1.050452749583598
$ws.Range("J23").Value = 1.060544826523528
$ws.Range("K23").Value = 1.064981600895847
$ws.Range("L23").Value = 1.07103217450944
$ws.Range("M23").Value = 1.077562705875741
$ws.Range("N23").Value = 1.062050921823382
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.055366406779507
$ws.Range("D24").Value = 1.062553630766937
$ws.Range("E24").Value = 1.06882388638867
$ws.Range("F24").Value = 1.075375116591255
$ws.Range("I24").Value = 1.050843729859307
$ws.Range("J24").Value = 1.061428491250354
$ws.Range("K24").Value = 1.065833571503122
$ws.Range("L24").Value = 1.07208326312599
$ws.Range("M24").Value = 1.078613288391882
$ws.Range("N24").Value = 1.06293584145548
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.056985864290414
$ws.Range("D25").Value = 1.063855018413424
$ws.Range("E25").Value = 1.070356880605381
$ws.Range("F25").Value = 1.076905191627113
$ws.Range("I25").Value = 1.051293730254057
$ws.Range("J25").Value = 1.062454058020433
$ws.Range("K25").Value = 1.066821912646208
$ws.Range("L25").Value = 1.07330462603446
$ws.Range("M25").Value = 1.079833905521387
$ws.Range("N25").Value = 1.063962864648006
